$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: tensorflow.blog post title/link update
$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝” 혼공 노트 증정 이벤트!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/01/20/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d-%ed%98%bc%ea%b3%b5-%eb%85%b8%ed%8a%b8-%ec%a6%9d%ec%a0%95-%ec%9d%b4%eb%b2%a4/"

# Row 23: Be the only one post title/link update
$ws.Range("D23").Value = "안녕하세요 딥러닝 논문읽기 모임 입니다!`n오늘 소개해 드릴 논문은 현재 많은 Image classification 분야에서 SOTA를 달성했던"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2672"

# Row 28: ropiens post title/link update
$ws.Range("D28").Value = "R-CNN : Region-based Convolutional Networks forAccurate Object Detection and Segmentation 리뷰 (공부 중)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/73"

# Row 37: dsba_seminar post title/link update
$ws.Range("D37").Value = "[Paper Review] Self-Supervised Learning by Cross-Modal Audio-Video Clustering"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1431&mod=document&pageid=1"

# Row 39: deadNstreet post title/link update
$ws.Range("D39").Value = "Machine Learning Wars: Deep Learning vs GBM"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Machine-Learning-Wars-Deep-Learning-vs-GBM-1"

# Row 40: insightCampus post title update (strip "New" notify span + its trailing tabs)
$ws.Range("D40").Value = "2020년 머신러닝 프로젝트 Top 10`t`t`t`t`t`t`t`t<span class=`"kboard-comments-count`"></span>"
